# Update "想去人数" (want-to-go count) values in column F across sheets
# 展览 (sheet1), 演出 (sheet2), and 全部类型 (sheet4).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsShow       = $wb.Worksheets.Item("演出")
$wsAll        = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$wsExhibition.Range("F2").Value  = 2965
$wsExhibition.Range("F6").Value  = 533
$wsExhibition.Range("F7").Value  = 66
$wsExhibition.Range("F9").Value  = 2929
$wsExhibition.Range("F12").Value = 7445
$wsExhibition.Range("F15").Value = 107
$wsExhibition.Range("F16").Value = 249
$wsExhibition.Range("F19").Value = 9096
$wsExhibition.Range("F27").Value = 113
$wsExhibition.Range("F30").Value = 91
$wsExhibition.Range("F31").Value = 67
$wsExhibition.Range("F33").Value = 2614
$wsExhibition.Range("F38").Value = 764
$wsExhibition.Range("F39").Value = 3895
$wsExhibition.Range("F40").Value = 205
$wsExhibition.Range("F41").Value = 35
$wsExhibition.Range("F43").Value = 78
$wsExhibition.Range("F46").Value = 11
$wsExhibition.Range("F47").Value = 58
$wsExhibition.Range("F48").Value = 30
$wsExhibition.Range("F49").Value = 56

# --- 演出 (sheet2) ---
$wsShow.Range("F7").Value = 123

# --- 全部类型 (sheet4) ---
$wsAll.Range("F3").Value  = 2965
$wsAll.Range("F8").Value  = 123
$wsAll.Range("F10").Value = 533
$wsAll.Range("F11").Value = 66
$wsAll.Range("F13").Value = 2929
$wsAll.Range("F18").Value = 7445
$wsAll.Range("F21").Value = 107
$wsAll.Range("F22").Value = 249
$wsAll.Range("F24").Value = 9096
$wsAll.Range("F30").Value = 113
$wsAll.Range("F31").Value = 91
$wsAll.Range("F32").Value = 67
$wsAll.Range("F34").Value = 2614
$wsAll.Range("F39").Value = 764
$wsAll.Range("F41").Value = 3895
$wsAll.Range("F42").Value = 205
$wsAll.Range("F43").Value = 35
$wsAll.Range("F47").Value = 58
$wsAll.Range("F48").Value = 30
$wsAll.Range("F49").Value = 56
